$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 655.7895
$ws.Cells.Item(15, 9).Value = 655.7895
$ws.Cells.Item(15, 11).Value = 1967.3685
$ws.Cells.Item(15, 13).Value = -1798.3685
$ws.Cells.Item(40, 8).Value = 3434.3333
$ws.Cells.Item(40, 10).Value = 5002
$ws.Cells.Item(40, 12).Value = 5002
$ws.Cells.Item(40, 14).Value = -5352
$ws.Cells.Item(113, 8).Value = 9068.714
$ws.Cells.Item(113, 9).Value = 6749.5
$ws.Cells.Item(113, 11).Value = 6749.5
$ws.Cells.Item(113, 13).Value = -3495.5
$ws.Cells.Item(115, 8).Value = 351.125
$ws.Cells.Item(115, 9).Value = 351.125
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 11).Value = 1053.375
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 13).Value = 513.625
$ws.Cells.Item(115, 14).ClearContents()
$ws.Cells.Item(127, 8).Value = 11683.857
$ws.Cells.Item(127, 9).Value = 800
$ws.Cells.Item(127, 11).Value = 2400
$ws.Cells.Item(127, 13).Value = 2560
$ws.Cells.Item(132, 8).Value = 4302.7427
$ws.Cells.Item(132, 9).Value = 4366.5757
$ws.Cells.Item(132, 11).Value = 13099.7271
$ws.Cells.Item(132, 13).Value = -10569.7271

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5646.8604
$ws.Cells.Item(32, 9).Value = 4551.075
$ws.Cells.Item(32, 11).Value = 4551.075
$ws.Cells.Item(32, 13).Value = -4264.075
$ws.Cells.Item(74, 8).Value = 2971.353
$ws.Cells.Item(74, 10).Value = 3239.8
$ws.Cells.Item(74, 12).Value = 3239.8
$ws.Cells.Item(74, 14).Value = -4987.8
$ws.Cells.Item(77, 8).Value = 2971.353
$ws.Cells.Item(77, 10).Value = 3239.8
$ws.Cells.Item(77, 12).Value = 16199
$ws.Cells.Item(77, 14).Value = -24935
$ws.Cells.Item(97, 8).Value = 4109.1904
$ws.Cells.Item(97, 9).Value = 4074.7222
$ws.Cells.Item(97, 11).Value = 4074.7222
$ws.Cells.Item(97, 13).Value = -3578.7222
$ws.Cells.Item(103, 8).Value = 252717.5
$ws.Cells.Item(103, 10).Value = 435435
$ws.Cells.Item(103, 12).Value = 435435
$ws.Cells.Item(103, 14).Value = -437779
$ws.Cells.Item(122, 8).Value = 41670730
$ws.Cells.Item(122, 9).Value = 5000
$ws.Cells.Item(122, 10).Value = 47622976
$ws.Cells.Item(122, 11).Value = 15000
$ws.Cells.Item(122, 12).Value = 142868928
$ws.Cells.Item(122, 13).Value = -12550
$ws.Cells.Item(122, 14).Value = -142873828
$ws.Cells.Item(132, 8).Value = 41865.188
$ws.Cells.Item(132, 9).Value = 4711.8076
$ws.Cells.Item(132, 11).Value = 14135.4228
$ws.Cells.Item(132, 13).Value = -11605.4228

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2800
$ws.Cells.Item(16, 9).Value = 2800
$ws.Cells.Item(16, 11).Value = 2800
$ws.Cells.Item(16, 13).Value = -2513
$ws.Cells.Item(31, 8).Value = 9618465
$ws.Cells.Item(31, 9).Value = 2319.75
$ws.Cells.Item(31, 10).Value = 41672284
$ws.Cells.Item(31, 11).Value = 2319.75
$ws.Cells.Item(31, 12).Value = 41672284
$ws.Cells.Item(31, 13).Value = -2024.75
$ws.Cells.Item(31, 14).Value = -41672874
$ws.Cells.Item(34, 8).Value = 9618465
$ws.Cells.Item(34, 9).Value = 2319.75
$ws.Cells.Item(34, 10).Value = 41672284
$ws.Cells.Item(34, 11).Value = 2319.75
$ws.Cells.Item(34, 12).Value = 41672284
$ws.Cells.Item(34, 13).Value = -2117.75
$ws.Cells.Item(34, 14).Value = -41672688
$ws.Cells.Item(58, 8).Value = 22919384
$ws.Cells.Item(58, 9).Value = 2499.8333
$ws.Cells.Item(58, 10).Value = 45836268
$ws.Cells.Item(58, 11).Value = 2499.8333
$ws.Cells.Item(58, 12).Value = 45836268
$ws.Cells.Item(58, 13).Value = -2296.8333
$ws.Cells.Item(58, 14).Value = -45836674
$ws.Cells.Item(99, 8).Value = 5166.3335
$ws.Cells.Item(99, 9).Value = 5166.3335
$ws.Cells.Item(99, 11).Value = 5166.3335
$ws.Cells.Item(99, 13).Value = -3668.3335
$ws.Cells.Item(105, 8).Value = 35670.332
$ws.Cells.Item(105, 10).Value = 35670.332
$ws.Cells.Item(105, 12).Value = 35670.332
$ws.Cells.Item(105, 14).Value = -39164.332
$ws.Cells.Item(113, 8).Value = 2800
$ws.Cells.Item(113, 9).Value = 2800
$ws.Cells.Item(113, 11).Value = 2800
$ws.Cells.Item(113, 13).Value = -630
$ws.Cells.Item(126, 8).Value = 5166.3335
$ws.Cells.Item(126, 9).Value = 5166.3335
$ws.Cells.Item(126, 11).Value = 15499.0005
$ws.Cells.Item(126, 13).Value = -13029.0005
$ws.Cells.Item(132, 8).Value = 3295.4243
$ws.Cells.Item(132, 9).Value = 3171.7
$ws.Cells.Item(132, 11).Value = 9515.099999999999
$ws.Cells.Item(132, 13).Value = -6985.099999999999
$ws.Cells.Item(134, 8).Value = 2626.394
$ws.Cells.Item(134, 9).Value = 2348.92
$ws.Cells.Item(134, 11).Value = 7046.76
$ws.Cells.Item(134, 13).Value = -4511.76
$ws.Cells.Item(136, 8).Value = 22919384
$ws.Cells.Item(136, 9).Value = 2499.8333
$ws.Cells.Item(136, 10).Value = 45836268
$ws.Cells.Item(136, 11).Value = 7499.499899999999
$ws.Cells.Item(136, 12).Value = 137508804
$ws.Cells.Item(136, 13).Value = -4949.499899999999
$ws.Cells.Item(136, 14).Value = -137513904

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 14).ClearContents()
$ws.Cells.Item(131, 8).Value = 45455830
$ws.Cells.Item(131, 10).Value = 1640.75
$ws.Cells.Item(131, 12).Value = 4922.25
$ws.Cells.Item(131, 14).Value = -15002.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 15217399
$ws.Cells.Item(11, 10).Value = 1819000
$ws.Cells.Item(11, 12).Value = 1819000
$ws.Cells.Item(11, 14).Value = -1819278
$ws.Cells.Item(80, 8).Value = 4991.2
$ws.Cells.Item(80, 9).Value = 3241.125
$ws.Cells.Item(80, 10).Value = 6991.2856
$ws.Cells.Item(80, 11).Value = 3241.125
$ws.Cells.Item(80, 12).Value = 6991.2856
$ws.Cells.Item(80, 13).Value = -2243.125
$ws.Cells.Item(80, 14).Value = -8987.2856
$ws.Cells.Item(83, 8).Value = 4991.2
$ws.Cells.Item(83, 9).Value = 3241.125
$ws.Cells.Item(83, 10).Value = 6991.2856
$ws.Cells.Item(83, 11).Value = 16205.625
$ws.Cells.Item(83, 12).Value = 34956.428
$ws.Cells.Item(83, 13).Value = -11213.625
$ws.Cells.Item(83, 14).Value = -44940.428
$ws.Cells.Item(132, 8).Value = 1667.6666
$ws.Cells.Item(132, 9).Value = 1400.7142
$ws.Cells.Item(132, 11).Value = 4202.142599999999
$ws.Cells.Item(132, 13).Value = -1672.142599999999
$ws.Cells.Item(136, 8).Value = 32317.857
$ws.Cells.Item(136, 10).Value = 32317.857
$ws.Cells.Item(136, 12).Value = 96953.571
$ws.Cells.Item(136, 14).Value = -102053.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(32, 8).Value = 5000
$ws.Cells.Item(32, 9).Value = 5000
$ws.Cells.Item(32, 11).Value = 5000
$ws.Cells.Item(32, 13).Value = -4683
$ws.Cells.Item(46, 8).Value = 3454.2727
$ws.Cells.Item(132, 8).Value = 1702.4375
$ws.Cells.Item(132, 9).Value = 1509.931
$ws.Cells.Item(132, 11).Value = 4529.793
$ws.Cells.Item(132, 13).Value = -1999.793
$ws.Cells.Item(133, 8).Value = 57720
$ws.Cells.Item(133, 10).Value = 57720
$ws.Cells.Item(133, 12).Value = 57720
$ws.Cells.Item(133, 14).Value = -62780
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 291669.84
$ws.Cells.Item(18, 9).Value = 667931.3
$ws.Cells.Item(18, 11).Value = 667931.3
$ws.Cells.Item(18, 13).Value = -667758.3
$ws.Cells.Item(62, 8).Value = 6966.1665
$ws.Cells.Item(65, 8).Value = 6966.1665
$ws.Cells.Item(95, 8).Value = 44998
$ws.Cells.Item(95, 10).Value = 44998
$ws.Cells.Item(95, 12).Value = 44998
$ws.Cells.Item(95, 14).Value = -50490
$ws.Cells.Item(96, 8).Value = 6672
$ws.Cells.Item(96, 9).Value = 8700
$ws.Cells.Item(96, 11).Value = 8700
$ws.Cells.Item(96, 13).Value = -7327
$ws.Cells.Item(100, 8).Value = 202001150
$ws.Cells.Item(100, 9).Value = 202001150
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 404002300
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).Value = -404001759
$ws.Cells.Item(100, 14).ClearContents()
